$wb = $excel.ActiveWorkbook

$wsGeneral = $wb.Worksheets.Item("General")
$wsGeneral.Rows.Item(384).Delete()

$wsAyko = $wb.Worksheets.Item("AYKO")
$wsAyko.Rows.Item(89).Delete()
